$d = $word.ActiveDocument

# --- Locate the paragraph that contains the "Website:" hyperlink ---
$findRng = $d.Content
$found = $findRng.Find.Execute("http://spotabee.buzz/", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$findRng.Expand(4)   # wdParagraph -> grow the hit to its whole paragraph
$hyperlinkParaEnd = $findRng.End

# --- 1. Remove the stray "_GoBack" bookmark that currently sits after the
#        "...back gardens" sentence ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Drop the trailing run of empty paragraphs that follow the hyperlink
#        paragraph, so that paragraph becomes the last one in the body
#        (directly followed by the section properties) ---
$endOfDoc = $d.Content.End
if ($endOfDoc -gt $hyperlinkParaEnd) {
    $trailing = $d.Range($hyperlinkParaEnd, $endOfDoc)
    $trailing.Delete()
}

# --- 3. Re-insert the "_GoBack" bookmark, now collapsed at the end of the
#        hyperlink paragraph (right after the URL, before the pilcrow) ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$paraEnd = $lastPara.Range.End - 1
$bmRange = $d.Range($paraEnd - 1, $paraEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
